$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered/top aligned) from existing header cell H1
# into the two new header cells I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row 2 data
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

# Row 3 data
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 8
